$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell C10 value from 18 to 100
$ws.Range("C10").Value = 100.0
